# Weekly update: insert a new price-report row for "Vega Modelo de Temuco - Alcachofa"
# at the top of the data block (row 50), pushing all subsequent rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 50. This shifts the existing rows 50..136 down to 51..137.
$ws.Rows.Item(50).Insert()

# Populate the newly inserted row 50 with the new weekly data point.
$ws.Cells.Item(50, 1).Value = 10
$ws.Cells.Item(50, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(50, 3).Value = "La Araucanía"
$ws.Cells.Item(50, 4).Value = 44482
$ws.Cells.Item(50, 5).Value = 9
$ws.Cells.Item(50, 6).Value = 100112013
$ws.Cells.Item(50, 7).Value = "Alcachofa"
$ws.Cells.Item(50, 8).Value = "Española"
$ws.Cells.Item(50, 9).Value = "Primera"
$ws.Cells.Item(50, 10).Value = 50
$ws.Cells.Item(50, 11).Value = 10000
$ws.Cells.Item(50, 12).Value = 10000
$ws.Cells.Item(50, 13).Value = 10000
$ws.Cells.Item(50, 14).Value = "$/caja 30 unidades"
$ws.Cells.Item(50, 15).Value = "Región Metropolitana"
$ws.Cells.Item(50, 16).Value = 333
$ws.Cells.Item(50, 17).Value = 30
$ws.Cells.Item(50, 18).Value = "Hortaliza"
